# Generate Report for Handback
# The 5ad28dad-3920-4125-bac1-dd1defa03890.md entry has been handed back
# (in sync with en-US) for both zh-cn and de-de locales. Update status,
# handback datetime, and clear the stale error detail message.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 5ad28dad... file; columns E (zh-cn) and F (de-de)
# hold the localization status text.
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is the 5ad28dad... file.
# Column C = Status, Column K = Latest Handback DateTime, Column P = Error Detail
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-31 15:02:51"
$zhcn.Range("P3").Value = ""

# de-de sheet: row 3 is the 5ad28dad... file.
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-31 15:03:00"
$dede.Range("P3").Value = ""
